$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4704.0415
$ws.Range("J88").Value = 6393.5884
$ws.Range("L88").Value = 6393.5884
$ws.Range("N88").Value = -7205.5884

$ws.Range("H91").Value = 4704.0415
$ws.Range("J91").Value = 6393.5884
$ws.Range("L91").Value = 6393.5884
$ws.Range("N91").Value = -9201.588400000001

$ws.Range("H118").Value = 2069.3333
$ws.Range("I118").Value = 690
$ws.Range("J118").Value = 2988.889
$ws.Range("K118").Value = 2070
$ws.Range("L118").Value = 8966.667000000001
$ws.Range("M118").Value = -413
$ws.Range("N118").Value = -12280.667

$ws.Range("H123").Value = 26800
$ws.Range("J123").Value = 26800
$ws.Range("L123").Value = 26800
$ws.Range("N123").Value = -36600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 984.4483
$ws.Range("I2").Value = 743.7083
$ws.Range("J2").Value = 2140
$ws.Range("K2").Value = 743.7083
$ws.Range("L2").Value = 2140
$ws.Range("M2").Value = -630.7083
$ws.Range("N2").Value = -2366

$ws.Range("H32").Value = 259213.11
$ws.Range("I32").Value = 298526.84
$ws.Range("K32").Value = 298526.84
$ws.Range("M32").Value = -298239.84

$ws.Range("H61").Value = 214434.23
$ws.Range("I61").Value = 1834.6
$ws.Range("J61").Value = 371915.44
$ws.Range("K61").Value = 1834.6
$ws.Range("L61").Value = 371915.44
$ws.Range("M61").Value = -1622.6
$ws.Range("N61").Value = -372339.44

$ws.Range("H116").Value = 984.4483
$ws.Range("I116").Value = 743.7083
$ws.Range("J116").Value = 2140
$ws.Range("K116").Value = 743.7083
$ws.Range("L116").Value = 2140
$ws.Range("M116").Value = 1550.2917
$ws.Range("N116").Value = -6728

$ws.Range("H132").Value = 7857.4546
$ws.Range("I132").Value = 5447.6
$ws.Range("J132").Value = 15388.25
$ws.Range("K132").Value = 16342.8
$ws.Range("L132").Value = 46164.75
$ws.Range("M132").Value = -13812.8
$ws.Range("N132").Value = -51224.75

$ws.Range("H136").Value = 214434.23
$ws.Range("I136").Value = 1834.6
$ws.Range("J136").Value = 371915.44
$ws.Range("K136").Value = 5503.799999999999
$ws.Range("L136").Value = 1115746.32
$ws.Range("M136").Value = -2953.799999999999
$ws.Range("N136").Value = -1120846.32

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 984.4483
$ws.Range("I3").Value = 743.7083
$ws.Range("J3").Value = 2140
$ws.Range("K3").Value = 743.7083
$ws.Range("L3").Value = 2140
$ws.Range("M3").Value = -629.7083
$ws.Range("N3").Value = -2368

$ws.Range("H20").Value = 27442.719
$ws.Range("I20").Value = 1537.56
$ws.Range("J20").Value = 73701.92999999999
$ws.Range("K20").Value = 1537.56
$ws.Range("L20").Value = 73701.92999999999
$ws.Range("M20").Value = -1290.56
$ws.Range("N20").Value = -74195.92999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16181.125
$ws.Range("I31").Value = 1145.7675
$ws.Range("J31").Value = 33654.65
$ws.Range("K31").Value = 1145.7675
$ws.Range("L31").Value = 33654.65
$ws.Range("M31").Value = -850.7674999999999
$ws.Range("N31").Value = -34244.65

$ws.Range("H34").Value = 16181.125
$ws.Range("I34").Value = 1145.7675
$ws.Range("J34").Value = 33654.65
$ws.Range("K34").Value = 1145.7675
$ws.Range("L34").Value = 33654.65
$ws.Range("M34").Value = -943.7674999999999
$ws.Range("N34").Value = -34058.65

$ws.Range("H86").Value = 3195.3215
$ws.Range("I86").Value = 2647.8235
$ws.Range("K86").Value = 2647.8235
$ws.Range("M86").Value = -1524.8235

$ws.Range("H89").Value = 3195.3215
$ws.Range("I89").Value = 2647.8235
$ws.Range("K89").Value = 13239.1175
$ws.Range("M89").Value = -7623.1175

$ws.Range("H107").Value = 708.44684
$ws.Range("I107").Value = 702.2381
$ws.Range("J107").Value = 760.6
$ws.Range("K107").Value = 702.2381
$ws.Range("L107").Value = 760.6
$ws.Range("M107").Value = 1217.7619
$ws.Range("N107").Value = -4600.6

$ws.Range("H124").Value = 18823.857
$ws.Range("J124").Value = 18823.857
$ws.Range("L124").Value = 18823.857
$ws.Range("N124").Value = -23733.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 617.9184
$ws.Range("I122").Value = 312
$ws.Range("J122").Value = 847.3570999999999
$ws.Range("K122").Value = 2808
$ws.Range("L122").Value = 7626.2139
$ws.Range("M122").Value = -358
$ws.Range("N122").Value = -12526.2139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3186
$ws.Range("I70").Value = 3572.2856
$ws.Range("J70").Value = 2799.7144
$ws.Range("K70").Value = 3572.2856
$ws.Range("L70").Value = 2799.7144
$ws.Range("M70").Value = -3302.2856
$ws.Range("N70").Value = -3339.7144

$ws.Range("H73").Value = 3186
$ws.Range("I73").Value = 3572.2856
$ws.Range("J73").Value = 2799.7144
$ws.Range("K73").Value = 3572.2856
$ws.Range("L73").Value = 2799.7144
$ws.Range("M73").Value = -2636.2856
$ws.Range("N73").Value = -4671.7144

$ws.Range("H132").Value = 21469.246
$ws.Range("I132").Value = 46321.48
$ws.Range("K132").Value = 138964.44
$ws.Range("M132").Value = -136434.44

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2144.125
$ws.Range("I7").Value = 2063.2727
$ws.Range("J7").Value = 2322
$ws.Range("K7").Value = 2063.2727
$ws.Range("L7").Value = 2322
$ws.Range("M7").Value = -1951.2727
$ws.Range("N7").Value = -2546

$ws.Range("H61").Value = 1743.3334
$ws.Range("I61").Value = 1405.2941
$ws.Range("J61").Value = 2185.3845
$ws.Range("K61").Value = 1405.2941
$ws.Range("L61").Value = 2185.3845
$ws.Range("M61").Value = -1203.2941
$ws.Range("N61").Value = -2589.3845

$ws.Range("H113").Value = 1743.3334
$ws.Range("I113").Value = 1405.2941
$ws.Range("J113").Value = 2185.3845
$ws.Range("K113").Value = 1405.2941
$ws.Range("L113").Value = 2185.3845
$ws.Range("M113").Value = 764.7058999999999
$ws.Range("N113").Value = -6525.3845

$ws.Range("H126").Value = 2144.125
$ws.Range("I126").Value = 2063.2727
$ws.Range("J126").Value = 2322
$ws.Range("K126").Value = 6189.8181
$ws.Range("L126").Value = 6966
$ws.Range("M126").Value = -3719.8181
$ws.Range("N126").Value = -11906

$ws.Range("H132").Value = 14933.185
$ws.Range("I132").Value = 8621.579
$ws.Range("J132").Value = 21244.79
$ws.Range("K132").Value = 25864.737
$ws.Range("L132").Value = 63734.37
$ws.Range("M132").Value = -23334.737
$ws.Range("N132").Value = -68794.37

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 9000
$ws.Range("J51").Value = 9000
$ws.Range("L51").Value = 9000
$ws.Range("N51").Value = -10020

$ws.Range("H107").Value = 477.66666
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 477.66666
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1432.99998
$ws.Range("N107").Value = -5272.999980000001
$ws.Range("M107").ClearContents()

$ws.Range("H122").Value = 33334168
$ws.Range("I122").Value = 40000640
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 120001920
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -119999470
$ws.Range("N122").Value = -10300
